$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 30 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(30, 1).Value = "Afmelding nieuwsbrief"
$logs.Cells.Item(30, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(30, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Cells.Item(30, 4).Value = "Afmelding"
$logs.Cells.Item(30, 6).Value = "2025-06-17 21:40:19"
$logs.Cells.Item(30, 7).Value = "Nee"

# Extend the conditional formatting ranges to include the new row
$catCond = $logs.Range("D2:D29").FormatConditions.Item(1)
$catCond.ModifyAppliesToRange($logs.Range("D2:D30"))

$answeredCond = $logs.Range("G2:G29").FormatConditions.Item(1)
$answeredCond.ModifyAppliesToRange($logs.Range("G2:G30"))

# --- Dashboard sheet: swap Afmelding / Bestelling rows and update count ---
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Cells.Item(4, 1).Value = "Afmelding"
$dashboard.Cells.Item(4, 2).Value = 4
$dashboard.Cells.Item(5, 1).Value = "Bestelling"
$dashboard.Cells.Item(5, 2).Value = 3
